# patRoon set-impl tracking sheet: add "getEICsForFGroups" entry and mark
# the fGroupsSet plotEIC method as done (commit: "plotEIC methods for
# fGroupsSet").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fGroups")

# Insert a new row above row 19 ("getFeatures"); this pushes getFeatures
# and everything below it down by one row, just like typing a new line
# in the middle of the tracking table in Excel.
$ws.Rows(19).Insert()

# Fill in the newly inserted row 19 with the new tracked method. Columns:
# B=as-is, C=almost as-is, D=implement, E=not supported, F=ionize, G=done
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# The plotEIC row (now shifted from row 33 down to row 34) is now
# implemented for fGroupsSet, so mark its "done" column too.
$ws.Range("G34").Value = "X"

# Leave the selection on the cell that was just finished, matching the
# workbook's last saved cursor position.
$ws.Range("G35").Select()
